# Test Case Fitur Topik - Membuat Materi
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Topik")

# --- Row 3 (new second "happy path" test case row, duplicate of row 2 but with multi-tag) ---
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Fundamental Katalon"
$ws.Range("C3").Value = "java-logo.jpg"
$ws.Range("D3").Value = "Fundamental Katalon yang mengenalkan java dari dasar."
$ws.Range("E3").Value = "JAVA,FUNDAMENTAL,INPUT"
$ws.Range("F3").Value = "Mudah"
$ws.Range("G3").Value = "Publik"

# --- Row 4 (missing thumbnail / oversized file test case) ---
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Fundamental Katalon"

# --- Row 5 (wrong file extension test case) ---
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "Fundamental Katalon"

# --- Row 6 (missing Judul/Topik test case) ---
$ws.Range("A6").Value = 5
$ws.Range("C6").Value = "java-logo.jpg"
$ws.Range("D6").Value = "Fundamental Katalon yang mengenalkan java dari dasar."
$ws.Range("E6").Value = "JAVA"
$ws.Range("F6").Value = "Mudah"
$ws.Range("G6").Value = "Publik"

# --- Row 7 (missing Thumbnail test case) ---
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "Fundamental Katalon"
$ws.Range("D7").Value = "Fundamental Katalon yang mengenalkan java dari dasar."
$ws.Range("E7").Value = "JAVA"
$ws.Range("F7").Value = "Mudah"
$ws.Range("G7").Value = "Publik"

# --- Row 8 (missing Deskripsi test case) ---
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "Fundamental Katalon"
$ws.Range("C8").Value = "java-logo.jpg"
$ws.Range("E8").Value = "JAVA"
$ws.Range("F8").Value = "Mudah"
$ws.Range("G8").Value = "Publik"

# --- Row 9 (missing Tag test case) ---
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "Fundamental Katalon"
$ws.Range("C9").Value = "java-logo.jpg"
$ws.Range("D9").Value = "Fundamental Katalon yang mengenalkan java dari dasar."
$ws.Range("F9").Value = "Mudah"
$ws.Range("G9").Value = "Publik"

# --- Row 10 (missing Tingkatan test case) ---
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "Fundamental Katalon"
$ws.Range("C10").Value = "java-logo.jpg"
$ws.Range("D10").Value = "Fundamental Katalon yang mengenalkan java dari dasar."
$ws.Range("E10").Value = "JAVA"
$ws.Range("G10").Value = "Publik"

# --- Row 11 (everything blank test case) ---
$ws.Range("A11").Value = 10

# --- Row 12 (missing Visibilitas test case) ---
$ws.Range("A12").Value = 11
$ws.Range("B12").Value = "Fundamental Katalon"
$ws.Range("C12").Value = "java-logo.jpg"
$ws.Range("D12").Value = "Fundamental Katalon yang mengenalkan java dari dasar."
$ws.Range("E12").Value = "JAVA"
$ws.Range("F12").Value = "Mudah"
$ws.Range("G12").Value = "Publik"

# --- Column J: test case IDs, filled top to bottom after the data rows ---
$ws.Range("J2").Value = "TC-Membuat materi-001"
$ws.Range("J3").Value = "TC-Membuat materi-002"
$ws.Range("J4").Value = "TC-Membuat materi-003"
$ws.Range("J5").Value = "TC-Membuat materi-004"
$ws.Range("J6").Value = "TC-Membuat materi-005"
$ws.Range("J7").Value = "TC-Membuat materi-006"
$ws.Range("J8").Value = "TC-Membuat materi-007"
$ws.Range("J9").Value = "TC-Membuat materi-008"
$ws.Range("J10").Value = "TC-Membuat materi-009"
$ws.Range("J11").Value = "TC-Membuat materi-010"
$ws.Range("J12").Value = "TC-Membuat materi-011"

# --- Edge-case descriptions added last in columns C4 / C5 ---
$ws.Range("C4").Value = "ukuran lebih dari 2 MB"
$ws.Range("C5").Value = "file extensi salah"

# --- Column widths / formatting ---
$ws.Columns.Item(3).ColumnWidth = 19.90625
$ws.Columns.Item(5).ColumnWidth = 24.08984375
$ws.Columns.Item(10).ColumnWidth = 21.54296875

# --- Selection matches the final edit position ---
$ws.Range("D12").Select()
